$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 3483153.2
$ws.Cells.Item(17, 10).Value = 3628022
$ws.Cells.Item(17, 12).Value = 10884066
$ws.Cells.Item(17, 14).Value = -10884402

$ws.Cells.Item(33, 8).Value = 8276.571
$ws.Cells.Item(33, 9).Value = 10422.2
$ws.Cells.Item(33, 11).Value = 10422.2
$ws.Cells.Item(33, 13).Value = -10193.2

$ws.Cells.Item(86, 8).Value = 6700.857
$ws.Cells.Item(86, 9).Value = 6001
$ws.Cells.Item(86, 11).Value = 6001
$ws.Cells.Item(86, 13).Value = -4878

$ws.Cells.Item(89, 8).Value = 6700.857
$ws.Cells.Item(89, 9).Value = 6001
$ws.Cells.Item(89, 11).Value = 30005
$ws.Cells.Item(89, 13).Value = -24389

$ws.Cells.Item(125, 8).Value = 1937
$ws.Cells.Item(125, 9).Value = 1374.5
$ws.Cells.Item(125, 11).Value = 12370.5
$ws.Cells.Item(125, 13).Value = -9910.5

$ws.Cells.Item(135, 8).Value = 1580.1305
$ws.Cells.Item(135, 9).Value = 1793.2632
$ws.Cells.Item(135, 11).Value = 16139.3688
$ws.Cells.Item(135, 13).Value = -13604.3688

$ws.Cells.Item(138, 8).Value = 2144.4827
$ws.Cells.Item(138, 9).Value = 1747.8889
$ws.Cells.Item(138, 10).Value = 7498.5
$ws.Cells.Item(138, 11).Value = 5243.6667
$ws.Cells.Item(138, 12).Value = 22495.5
$ws.Cells.Item(138, 13).Value = -103.6666999999998
$ws.Cells.Item(138, 14).Value = -32775.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3722.383
$ws.Cells.Item(61, 9).Value = 2807.7693
$ws.Cells.Item(61, 10).Value = 4854.7617
$ws.Cells.Item(61, 11).Value = 2807.7693
$ws.Cells.Item(61, 12).Value = 4854.7617
$ws.Cells.Item(61, 13).Value = -2595.7693
$ws.Cells.Item(61, 14).Value = -5278.7617

$ws.Cells.Item(97, 8).Value = 607.04877
$ws.Cells.Item(97, 9).Value = 649.94446
$ws.Cells.Item(97, 11).Value = 649.94446
$ws.Cells.Item(97, 13).Value = -153.94446

$ws.Cells.Item(122, 8).Value = 1896.6666
$ws.Cells.Item(122, 9).Value = 276
$ws.Cells.Item(122, 11).Value = 828
$ws.Cells.Item(122, 13).Value = 1622

$ws.Cells.Item(132, 8).Value = 10477.75
$ws.Cells.Item(132, 9).Value = 10890.833
$ws.Cells.Item(132, 11).Value = 32672.499
$ws.Cells.Item(132, 13).Value = -30142.499

$ws.Cells.Item(136, 8).Value = 3722.383
$ws.Cells.Item(136, 9).Value = 2807.7693
$ws.Cells.Item(136, 10).Value = 4854.7617
$ws.Cells.Item(136, 11).Value = 8423.3079
$ws.Cells.Item(136, 12).Value = 14564.2851
$ws.Cells.Item(136, 13).Value = -5873.3079
$ws.Cells.Item(136, 14).Value = -19664.2851

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3166
$ws.Cells.Item(20, 10).Value = 1999
$ws.Cells.Item(20, 12).Value = 1999
$ws.Cells.Item(20, 14).Value = -2493

$ws.Cells.Item(94, 8).Value = 490.18518
$ws.Cells.Item(94, 9).Value = 468.66666
$ws.Cells.Item(94, 11).Value = 468.66666
$ws.Cells.Item(94, 13).Value = -17.66665999999998

$ws.Cells.Item(107, 8).Value = 1174.5416
$ws.Cells.Item(107, 9).Value = 1061
$ws.Cells.Item(107, 10).Value = 1969.3334
$ws.Cells.Item(107, 11).Value = 1061
$ws.Cells.Item(107, 12).Value = 1969.3334
$ws.Cells.Item(107, 13).Value = 859
$ws.Cells.Item(107, 14).Value = -5809.3334

$ws.Cells.Item(134, 8).Value = 14683.923
$ws.Cells.Item(134, 9).Value = 12987.4375
$ws.Cells.Item(134, 11).Value = 38962.3125
$ws.Cells.Item(134, 13).Value = -36427.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3464.9312
$ws.Cells.Item(31, 9).Value = 1640.1538
$ws.Cells.Item(31, 10).Value = 4947.5625
$ws.Cells.Item(31, 11).Value = 1640.1538
$ws.Cells.Item(31, 12).Value = 4947.5625
$ws.Cells.Item(31, 13).Value = -1345.1538
$ws.Cells.Item(31, 14).Value = -5537.5625

$ws.Cells.Item(34, 8).Value = 3464.9312
$ws.Cells.Item(34, 9).Value = 1640.1538
$ws.Cells.Item(34, 10).Value = 4947.5625
$ws.Cells.Item(34, 11).Value = 1640.1538
$ws.Cells.Item(34, 12).Value = 4947.5625
$ws.Cells.Item(34, 13).Value = -1438.1538
$ws.Cells.Item(34, 14).Value = -5351.5625

$ws.Cells.Item(70, 8).Value = 99988.75
$ws.Cells.Item(70, 10).Value = 99988.75
$ws.Cells.Item(70, 12).Value = 99988.75
$ws.Cells.Item(70, 14).Value = -100618.75

$ws.Cells.Item(73, 8).Value = 99988.75
$ws.Cells.Item(73, 10).Value = 99988.75
$ws.Cells.Item(73, 12).Value = 99988.75
$ws.Cells.Item(73, 14).Value = -102172.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = $null
$ws.Cells.Item(63, 14).Value = $null

$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 13).Value = $null
$ws.Cells.Item(66, 14).Value = $null

$ws.Cells.Item(80, 8).Value = 5999.5
$ws.Cells.Item(80, 10).Value = 5999.5
$ws.Cells.Item(80, 12).Value = 17998.5
$ws.Cells.Item(80, 14).Value = -19870.5

$ws.Cells.Item(81, 8).Value = 14856.379
$ws.Cells.Item(81, 10).Value = 17280.912
$ws.Cells.Item(81, 12).Value = 51842.736
$ws.Cells.Item(81, 14).Value = -54088.736

$ws.Cells.Item(83, 8).Value = 5999.5
$ws.Cells.Item(83, 10).Value = 5999.5
$ws.Cells.Item(83, 12).Value = 53995.5
$ws.Cells.Item(83, 14).Value = -63355.5

$ws.Cells.Item(84, 8).Value = 14856.379
$ws.Cells.Item(84, 10).Value = 17280.912
$ws.Cells.Item(84, 12).Value = 155528.208
$ws.Cells.Item(84, 14).Value = -166760.208

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 16371.667
$ws.Cells.Item(15, 10).Value = 16371.667
$ws.Cells.Item(15, 12).Value = 16371.667
$ws.Cells.Item(15, 14).Value = -16947.667

$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).Value = $null
$ws.Cells.Item(80, 14).Value = $null

$ws.Cells.Item(81, 8).Value = 16371.667
$ws.Cells.Item(81, 10).Value = 16371.667
$ws.Cells.Item(81, 12).Value = 16371.667
$ws.Cells.Item(81, 14).Value = -18367.667

$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).Value = $null
$ws.Cells.Item(83, 14).Value = $null

$ws.Cells.Item(84, 8).Value = 16371.667
$ws.Cells.Item(84, 10).Value = 16371.667
$ws.Cells.Item(84, 12).Value = 49115.001
$ws.Cells.Item(84, 14).Value = -59099.001

$ws.Cells.Item(132, 8).Value = 3011.8235
$ws.Cells.Item(132, 9).Value = 2993.8262
$ws.Cells.Item(132, 11).Value = 8981.4786
$ws.Cells.Item(132, 13).Value = -6451.4786

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 2754.5
$ws.Cells.Item(4, 9).Value = 4009
$ws.Cells.Item(4, 10).Value = 1500
$ws.Cells.Item(4, 11).Value = 4009
$ws.Cells.Item(4, 12).Value = 1500
$ws.Cells.Item(4, 13).Value = -3896
$ws.Cells.Item(4, 14).Value = -1726

$ws.Cells.Item(16, 8).Value = 1628.25
$ws.Cells.Item(16, 9).Value = 1902.7
$ws.Cells.Item(16, 10).Value = 256
$ws.Cells.Item(16, 11).Value = 1902.7
$ws.Cells.Item(16, 12).Value = 256
$ws.Cells.Item(16, 13).Value = -1732.7
$ws.Cells.Item(16, 14).Value = -596

$ws.Cells.Item(22, 8).Value = 2545.6
$ws.Cells.Item(22, 10).Value = 2019.5
$ws.Cells.Item(22, 12).Value = 2019.5
$ws.Cells.Item(22, 14).Value = -2609.5

$ws.Cells.Item(27, 8).Value = 2545.6
$ws.Cells.Item(27, 10).Value = 2019.5
$ws.Cells.Item(27, 12).Value = 2019.5
$ws.Cells.Item(27, 14).Value = -2233.5

$ws.Cells.Item(28, 8).Value = 2754.5
$ws.Cells.Item(28, 9).Value = 4009
$ws.Cells.Item(28, 10).Value = 1500
$ws.Cells.Item(28, 11).Value = 4009
$ws.Cells.Item(28, 12).Value = 1500
$ws.Cells.Item(28, 13).Value = -3777
$ws.Cells.Item(28, 14).Value = -1964

$ws.Cells.Item(37, 8).Value = 2754.5
$ws.Cells.Item(37, 9).Value = 4009
$ws.Cells.Item(37, 10).Value = 1500
$ws.Cells.Item(37, 11).Value = 4009
$ws.Cells.Item(37, 12).Value = 1500
$ws.Cells.Item(37, 13).Value = -3902
$ws.Cells.Item(37, 14).Value = -1714

$ws.Cells.Item(40, 8).Value = 2500
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 13).Value = $null

$ws.Cells.Item(61, 8).Value = 2128.5557
$ws.Cells.Item(61, 9).Value = 1996.8
$ws.Cells.Item(61, 11).Value = 1996.8
$ws.Cells.Item(61, 13).Value = -1794.8

$ws.Cells.Item(82, 8).Value = 1407.5
$ws.Cells.Item(82, 10).Value = 1682.3334
$ws.Cells.Item(82, 12).Value = 1682.3334
$ws.Cells.Item(82, 14).Value = -2404.3334

$ws.Cells.Item(85, 8).Value = 1407.5
$ws.Cells.Item(85, 10).Value = 1682.3334
$ws.Cells.Item(85, 12).Value = 1682.3334
$ws.Cells.Item(85, 14).Value = -4178.3334

$ws.Cells.Item(113, 8).Value = 2128.5557
$ws.Cells.Item(113, 9).Value = 1996.8
$ws.Cells.Item(113, 11).Value = 1996.8
$ws.Cells.Item(113, 13).Value = 173.2

$ws.Cells.Item(122, 8).Value = 3661.182
$ws.Cells.Item(122, 9).Value = 2247.5
$ws.Cells.Item(122, 10).Value = 5357.6
$ws.Cells.Item(122, 11).Value = 6742.5
$ws.Cells.Item(122, 12).Value = 16072.8
$ws.Cells.Item(122, 13).Value = -4292.5
$ws.Cells.Item(122, 14).Value = -20972.8

$ws.Cells.Item(136, 8).Value = 4276482
$ws.Cells.Item(136, 9).Value = 5293529
$ws.Cells.Item(136, 11).Value = 15880587
$ws.Cells.Item(136, 13).Value = -15878037

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 990.2857
$ws.Cells.Item(100, 9).Value = 870.4167
$ws.Cells.Item(100, 11).Value = 1740.8334
$ws.Cells.Item(100, 13).Value = -1199.8334

$ws.Cells.Item(107, 8).Value = 1773.92
$ws.Cells.Item(107, 10).Value = 1657.3572
$ws.Cells.Item(107, 12).Value = 4972.071599999999
$ws.Cells.Item(107, 14).Value = -8812.071599999999
